$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 6231.33
$ws1.Range("L3").Value = 537.34
$ws1.Range("D28").Value = 1831.68
$ws1.Range("M92").Value = 2521.53
$ws1.Range("M93").Value = 3887.15
$ws1.Range("M109").Value = 1054.31
$ws1.Range("L110").Value = 2780.75
$ws1.Range("M110").Value = 4359.41
$ws1.Range("H115").Value = 1161
$ws1.Range("I115").Value = 43.2
$ws1.Range("I122").Value = 626.4
$ws1.Range("P122").Value = 316.28
$ws1.Range("L134").Value = 1140.48
$ws1.Range("M134").Value = 3999.25
$ws1.Range("C135").Value = 518.4
$ws1.Range("M151").Value = 2272.64
$ws1.Range("D152").Value = 190.08
$ws1.Range("H152").Value = 1698.3
$ws1.Range("G154").Value = 166.43
$ws1.Range("M154").Value = 13047.51
$ws1.Range("L158").Value = 525.2
$ws1.Range("D163").Value = 1419.8
$ws1.Range("H163").Value = 2293.19
$ws1.Range("I163").Value = 255.6

# row 334 counter labels ("N de 332")
$ws1.Range("C334").Value = "7 de 332"
$ws1.Range("D334").Value = "28 de 332"
$ws1.Range("G334").Value = "1 de 332"
$ws1.Range("H334").Value = "17 de 332"
$ws1.Range("I334").Value = "21 de 332"
$ws1.Range("L334").Value = "39 de 332"
$ws1.Range("M334").Value = "72 de 332"
$ws1.Range("P334").Value = "4 de 332"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 6231.33
$ws2.Range("F3").Value = 687.03
$ws2.Range("F28").Value = 1831.68
$ws2.Range("F92").Value = 6849.29
$ws2.Range("F93").Value = 6679.65
$ws2.Range("F109").Value = 1054.31
$ws2.Range("F110").Value = 8432.22
$ws2.Range("F115").Value = 3417.52
$ws2.Range("F122").Value = 11589.38
$ws2.Range("F134").Value = 6441.13
$ws2.Range("F135").Value = 518.4
$ws2.Range("F155").Value = 2272.64
$ws2.Range("F156").Value = 6897.82
$ws2.Range("F158").Value = 14573.17
$ws2.Range("F162").Value = 6213.78
$ws2.Range("F167").Value = 10848.32
$ws2.Range("F338").Value = 369160.07

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 3592.51
$ws3.Range("E3").Value = 5242.06354940916
$ws3.Range("F3").Value = 0.4066421519848302
$ws3.Range("D11").Value = 4298.72
$ws3.Range("E11").Value = -1376.49541814726
$ws3.Range("F11").Value = 1.471043679084562
$ws3.Range("D12").Value = 20190.39
$ws3.Range("E12").Value = 2243.3653751766
$ws3.Range("F12").Value = 0.9000004529933081
$ws3.Range("D28").Value = 1321.92
$ws3.Range("E28").Value = 4875.66402943659
$ws3.Range("F28").Value = 0.2132960188552979
$ws3.Range("D32").Value = 4072.5
$ws3.Range("E32").Value = -1164.91631853974
$ws3.Range("F32").Value = 1.400647563806209
$ws3.Range("D33").Value = 2753.1
$ws3.Range("E33").Value = -1866.388983712426
$ws3.Range("F33").Value = 3.104844700730691
$ws3.Range("D34").Value = 316.28
$ws3.Range("E34").Value = 1030.12488751609
$ws3.Range("F34").Value = 0.2349070498276993
$ws3.Range("D37").Value = 15667.05
$ws3.Range("E37").Value = 2164.364398465401
$ws3.Range("F37").Value = 0.878620711172992
$ws3.Range("D38").Value = 51474.1
$ws3.Range("E38").Value = 10389.6203947566
$ws3.Range("F38").Value = 0.8320563275461008
$ws3.Range("D42").Value = 2983.64
$ws3.Range("E42").Value = 2520.97890386263
$ws3.Range("F42").Value = 0.5420248071862629
$ws3.Range("D44").Value = 166.43
$ws3.Range("E44").Value = -16.43000000000001
$ws3.Range("F44").Value = 1.109533333333333
$ws3.Range("D45").Value = 5690.69
$ws3.Range("E45").Value = -2783.10631853974
$ws3.Range("F45").Value = 1.957188725568165
$ws3.Range("D46").Value = 383.4
$ws3.Range("E46").Value = 503.311016287574
$ws3.Range("F46").Value = 0.4323843878755392
$ws3.Range("D50").Value = 1917.23
$ws3.Range("E50").Value = 3927.21916370549
$ws3.Range("F50").Value = 0.3280428910060774
$ws3.Range("D51").Value = 48945.04
$ws3.Range("E51").Value = -12121.3969078829
$ws3.Range("F51").Value = 1.329174299174048
$ws3.Range("D81").Value = 12588.88
$ws3.Range("E81").Value = 7411.120000000001
$ws3.Range("F81").Value = 0.629444
$ws3.Range("D97").Value = 395471.6699999999
$ws3.Range("E97").Value = 100043.9206021116
$ws3.Range("F97").Value = 0.7981013665371335
